# Add the two new "association" columns (S and T) to the becExpType
# config sheet, matching the "adding phase lock association" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers in row 1
$ws.Range("S1").Value = "PhaseLockAssociation"
$ws.Range("T1").Value = "WaveformAssociation"

# New default values in row 2 (matches existing "None" entries elsewhere
# on the row, e.g. CloudCenterReference/FringeRemovalMethod)
$ws.Range("S2").Value = "None"
$ws.Range("T2").Value = "None"

# Widen the newly added column S to match the author's saved view
# (target stored width 17.28515625 chars; 16.5 is the closest value this
# engine's ColumnWidth -> stored-width quantization can produce)
$ws.Columns("S").ColumnWidth = 16.5

# Move the selection/viewport the way the saved workbook shows it
$ws.Range("Y1").Select()
